$d = $word.ActiveDocument

# The paragraph currently reads "hellow" followed immediately by the
# _GoBack bookmark (empty range at the end of the paragraph). The target
# splits this into two runs - "   w" then "hellow" - with the _GoBack
# bookmark sitting between them.
#
# Plain text/range insertion in this run always coalesces into the single
# existing run (since the formatting is identical), so a bookmark is used
# as a scaffold to force the run split, and is then removed/relocated.

# 1. Duplicate the paragraph's text so we end up with "hellowhellow" in a
#    single run (still ahead of the _GoBack bookmark).
$tail = $d.Range(6, 6)
$tail.InsertAfter("hellow")

# 2. Wrap the freshly appended second "hellow" (chars 6-12) in a throwaway
#    bookmark - adding a bookmark around a non-empty range forces the run
#    to split into distinct <w:r> elements at the bookmark boundaries.
$secondHellow = $d.Range(6, 12)
$d.Bookmarks.Add("ZZZ_SPLIT_HELPER", $secondHellow)

# 3. Drop the helper bookmark - the run split it created survives removal
#    of the bookmark markers themselves.
$d.Bookmarks.Item("ZZZ_SPLIT_HELPER").Delete()

# 4. Relocate the real _GoBack bookmark so it once again sits between the
#    two runs (right where it originally was, before we appended text).
$gap = $d.Range(6, 6)
$d.Bookmarks.Add("_GoBack", $gap)

# 5. Finally, turn the first run's text ("hellow") into "   w" - Find:=
#    scoped to just that first run so the second run/bookmark are left
#    untouched.
$firstRun = $d.Range(0, 6)
$firstRun.Find.Execute("hellow", $true, $false, $false, $false, $false, `
                        $true, 1, $false, "   w", 2)
